# Insert a new weekly price record at row 230 (Hortaliza, Feria Lagunitas de
# Puerto Montt - Cebollín), pushing the existing row 230 (and every row
# after it) down by one. This grows the used range from A1:R354 to A1:R355.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 230..354 down to 231..355, leaving a blank row 230 to fill in.
$ws.Rows("230:230").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A230").Value = 4
$ws.Range("B230").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C230").Value = "Los Lagos"
$ws.Range("D230").Value = 44873
$ws.Range("E230").Value = 10
$ws.Range("F230").Value = 100112037
$ws.Range("G230").Value = "Cebollín"
$ws.Range("H230").Value = "Sin especificar"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 180
$ws.Range("K230").Value = 6500
$ws.Range("L230").Value = 6500
$ws.Range("M230").Value = 6500
$ws.Range("N230").Value = "`$/paquete 36 unidades"
$ws.Range("O230").Value = "Región Metropolitana"
$ws.Range("P230").Value = 181
$ws.Range("Q230").Value = 36
$ws.Range("R230").Value = "Hortaliza"
